$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = 0.9223046214701632
$ws.Range("J2").Value = 0.9223046214701632
$ws.Range("M2").Value = 29.52617166666667
$ws.Range("N2").Value = 88.57851500000001
$ws.Range("O2").Value = 0.3218391660320701
$ws.Range("P2").Value = 0.3218391660320701
$ws.Range("Q2").Value = 49.26503747543834
$ws.Range("R2").Value = 443.3853372789451
$ws.Range("S2").Value = 0.2968337502014814
$ws.Range("T2").Value = 0.2968337502014814
$ws.Range("I3").Value = 0.9223046214701632
$ws.Range("J3").Value = 0.9223046214701632
$ws.Range("O3").Value = 0.4328989896002822
$ws.Range("P3").Value = 0.4328989896002822
$ws.Range("S3").Value = 0.3992647387381044
$ws.Range("T3").Value = 0.3992647387381044
$ws.Range("I4").Value = 0.9223046214701632
$ws.Range("J4").Value = 0.9223046214701632
$ws.Range("M4").Value = 22.50081433333333
$ws.Range("N4").Value = 67.502443
$ws.Range("O4").Value = 0.2452618443676477
$ws.Range("P4").Value = 0.2452618443676476
$ws.Range("Q4").Value = 37.54308123226767
$ws.Range("R4").Value = 337.887731090409
$ws.Range("S4").Value = 0.2262061325305774
$ws.Range("T4").Value = 0.2262061325305773
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.140557
$ws.Range("H5").Value = 0.421671
$ws.Range("I5").Value = 0.07769537852983674
$ws.Range("J5").Value = 0.07769537852983674
$ws.Range("M5").Value = 29.52617166666667
$ws.Range("N5").Value = 88.57851500000001
$ws.Range("O5").Value = 0.3218391660320701
$ws.Range("P5").Value = 0.3218391660320701
$ws.Range("Q5").Value = 4.150110110951667
$ws.Range("R5").Value = 37.350990998565
$ws.Range("S5").Value = 0.02500541583058866
$ws.Range("T5").Value = 0.02500541583058866
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.140557
$ws.Range("H6").Value = 0.421671
$ws.Range("I6").Value = 0.07769537852983674
$ws.Range("J6").Value = 0.07769537852983674
$ws.Range("O6").Value = 0.4328989896002822
$ws.Range("P6").Value = 0.4328989896002822
$ws.Range("Q6").Value = 5.582224487811001
$ws.Range("R6").Value = 50.240020390299
$ws.Range("S6").Value = 0.03363425086217778
$ws.Range("T6").Value = 0.03363425086217778
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.140557
$ws.Range("H7").Value = 0.421671
$ws.Range("I7").Value = 0.07769537852983674
$ws.Range("J7").Value = 0.07769537852983674
$ws.Range("M7").Value = 22.50081433333333
$ws.Range("N7").Value = 67.502443
$ws.Range("O7").Value = 0.2452618443676477
$ws.Range("P7").Value = 0.2452618443676476
$ws.Range("Q7").Value = 3.162646960250334
$ws.Range("R7").Value = 28.463822642253
$ws.Range("S7").Value = 0.01905571183707029
$ws.Range("T7").Value = 0.01905571183707029
